# run_all bash file
# - drop the short-lived "Properties" sheet (tab2), keep "Properties_more"
#   (tab1) but rename it to "Properties" so it becomes the sole sheet.
# - replace its filler data (1000..30000) with a geometric "doubling-step"
#   series driven by formulas: A2 = 16000*1/25, then each subsequent cell
#   adds $A$2, down to row 26 (value 16000).
# - leave the active cell on A26 / the sole tab selected.

$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

# 1) Remove the second sheet ("Properties") - its two sample rows are gone.
$wb.Worksheets.Item("Properties").Delete() | Out-Null

# 2) The remaining sheet ("Properties_more") becomes the new "Properties".
$ws = $wb.Worksheets.Item("Properties_more")
$ws.Name = "Properties"

# 3) Wipe the old filler values (rows 2-31) before laying down new formulas.
$ws.Range("A2:A31").ClearContents()

# 4) Lay down the new formula series in A2:A26.
$ws.Range("A2").Formula = "=16000*1/25"
$ws.Range("A3").Formula = "=A2+`$A`$2"
for ($r = 4; $r -le 26; $r++) {
    $prev = $r - 1
    $ws.Range("A$r").Formula = "=A$prev+`$A`$2"
}

# 5) Match the saved view state: sole tab selected, active cell on the
#    last populated row.
$ws.Activate()
$ws.Range("A26").Select() | Out-Null
